# Update dados ADD (faturamento diario) - "atualizei dados da bibi e add"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Direct value corrections for a few existing August rows (B column = total_venda)
$ws.Range("B3").Value = 30358.01
$ws.Range("B8").Value = 12549.55
$ws.Range("B9").Value = 115299.64
$ws.Range("B12").Value = 50257.92999999999

# 2) Insert a new daily record (dia 18, agosto/2025) right after the existing
#    August rows, pushing every subsequent row (July/June/May data) down by one.
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 18
$ws.Range("B13").Value = 9371.9
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "08/2025"
